$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header labels for the two added rows (will become new shared strings)
$ws.Range("A12").Value = "Temp (F)"
$ws.Range("B12").Value = 70

$ws.Range("A13").Value = "Pressure  (kPa)"
$ws.Range("B13").Value = 100.25

# Column A width auto-fit (bestFit) similar to column D already present
$ws.Columns.Item(1).ColumnWidth = 14.28515625

$ws.Range("B16").Select()
